# lab2/assignment1/martin/data/samples.xlsx - update column F (Y) values
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = 4
$ws.Range("F3").Value = 8
$ws.Range("F4").Value = 12
$ws.Range("F5").Value = 16
$ws.Range("F6").Value = 20
$ws.Range("F7").Value = 24
$ws.Range("F8").Value = 28

# Move the active selection to F11 (as recorded on save in the source file)
[void]$ws.Range("F11").Select()
